# Apply "Week 15" simulation updates to the Players Data workbook.
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$wsRushing = $wb.Worksheets.Item("Rushing")

$wsRushing.Range("C2").Value = 11
$wsRushing.Range("D2").Value = 13
$wsRushing.Range("E2").Value = 7
$wsRushing.Range("F2").Value = 10

$wsRushing.Range("C3").Value = 65
$wsRushing.Range("D3").Value = 34
$wsRushing.Range("F3").Value = 12

$wsRushing.Range("C4").Value = 28
$wsRushing.Range("F4").Value = 9

$wsRushing.Range("C6").Value = 8
$wsRushing.Range("D6").Value = 5
$wsRushing.Range("E6").Value = 1

$wsRushing.Range("E7").Value = 2
$wsRushing.Range("F7").Value = 1

$wsRushing.Range("C9").Value = 5

# --- Receiving sheet ---
$wsReceiving = $wb.Worksheets.Item("Receiving")

$wsReceiving.Range("C2").Value = 33
$wsReceiving.Range("D2").Value = 27

$wsReceiving.Range("C3").Value = 26
$wsReceiving.Range("D3").Value = 19
$wsReceiving.Range("E3").Value = 4
$wsReceiving.Range("F3").Value = 4
$wsReceiving.Range("G3").Value = 4
$wsReceiving.Range("H3").Value = 3

$wsReceiving.Range("C5").Value = 3
$wsReceiving.Range("D5").Value = 3

$wsReceiving.Range("C6").Value = 101
$wsReceiving.Range("D6").Value = 79
$wsReceiving.Range("E6").Value = 29
$wsReceiving.Range("F6").Value = 10

$wsReceiving.Range("C7").Value = 50
$wsReceiving.Range("D7").Value = 39
$wsReceiving.Range("E7").Value = 12
$wsReceiving.Range("F7").Value = 5

$wsReceiving.Range("E8").Value = 13
$wsReceiving.Range("F8").Value = 9

$wsReceiving.Range("C11").Value = 6
$wsReceiving.Range("D11").Value = 3
$wsReceiving.Range("G11").Value = 3
$wsReceiving.Range("H11").Value = 1

$wsReceiving.Range("C12").Value = 89
$wsReceiving.Range("D12").Value = 62

$wb.Save()
